# Finish updating extraction data
# Appends the tail rows (61-69) of task-log entries to Sheet1, matching the
# "Finish updating extraction data" commit: new RQ2.1/RQ2.2/RQ2.3/RQ3 task
# rows with their time-spent (minutes) values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row, Task text (col A), Count (col B), Minutes (col C)
$rows = @(
    @(61, "New values of RQ2.1 on deployment, add new values", 1, 30),
    @(62, "Revampe values for RQ2.2 on on-chain data", 1, 40),
    @(63, "Revamp and add values to RQ2.2", 1, 65),
    @(64, "Make figures for RQ2.2", 1, 10),
    @(65, "Revamp values for RQ2.2 on on-chain logic", 1, 35),
    @(66, "Add values of new papers to RQ2.2 on on-chain logic", 1, 35),
    @(67, "Draw new figures for RQ2.2 on on-chain logic", 1, 5),
    @(68, "Update and add values for RQ2.3 ", 1, 45),
    @(69, "Update values for RQ3, add new values, draw figures", 1, 35)
)

foreach ($r in $rows) {
    $rowNum = $r[0]
    $ws.Cells.Item($rowNum, 1).Value = $r[1]
    $ws.Cells.Item($rowNum, 2).Value = $r[2]
    $ws.Cells.Item($rowNum, 3).Value = $r[3]
}

# Match the saved selection state: the bottom block C61:C69 ends up selected
# with the last cell as the active one.
$ws.Range("C61:C69").Select()
